$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("F2").Value = "Persiapan Juziyah 1 (Juziyah tanggal 6) - Orangnya sakit"

# --- Row 3 ---
$ws.Range("D3").Value = "3 kali (juz 8), 3 kali (juz 7)"
$ws.Range("E3").Value = "10 juz"
$ws.Range("F3").Value = ""

# --- Row 4 ---
$ws.Range("D4").Value = "(Di juz 2, sampai 3A)"
$ws.Range("E4").Value = "5 juz"
$ws.Range("F4").Value = "Persiapan Juziyah 2 (Juziyah tanggal 23)"

# --- Row 6 ---
$ws.Range("E6").Value = "9 juz (juz 8 kurang lancar)"
$ws.Range("F6").Value = "Persiapan Syahadah (udah siap juz 5, 6, 7, dan 9) Syahadah pertengahan Desember"

# --- Row 7 ---
$ws.Range("B7").Value = "Juz 4, 6A"
$ws.Range("C7").Value = "5 halaman"
$ws.Range("D7").Value = "5 kali"
$ws.Range("E7").Value = "5 juz"

# --- Row 8 ---
$ws.Range("D8").Value = "(Di juz 3, sampai 4B)"
$ws.Range("E8").Value = "4 juz (kurang juz 3)"

# --- Row 9 ---
$ws.Range("E9").Value = "9 juz (juz 9 kurang setengah)"
$ws.Range("F9").Value = "Persiapan Syahadah (udah siap juz 5, 6, 7, dan 8) Syahadah akhir Desember"

# --- Row 10 ---
$ws.Range("D10").Value = "(Di juz 5, sampai 2A)"
$ws.Range("E10").Value = "5 juz (juz 2 dikit lagi)"
$ws.Range("F10").Value = "Persiapan Juziyah 5 (Juziyah tanggal 18)"

# --- Row 12 (new) ---
$ws.Range("A12").Value = "Radja"
$ws.Range("B12").Value = "Juz 1, 9A"
$ws.Range("C12").Value = "8 halaman"
$ws.Range("D12").Value = "6 kali"
$ws.Range("E12").Value = "2 juz"

# --- Row 13 (new) ---
$ws.Range("A13").Value = "Faiq"
$ws.Range("B13").Value = "Juz 29, 1A"
$ws.Range("C13").Value = "7 halaman"
$ws.Range("D13").Value = "5 kali"
$ws.Range("E13").Value = "2 juz"

# --- Row 14 (new) ---
$ws.Range("A14").Value = "Shinra"
$ws.Range("B14").Value = "Juz 1, 10B"
$ws.Range("C14").Value = "20 halaman"
$ws.Range("D14").Value = "6 kali"
$ws.Range("E14").Value = "2 juz"
$ws.Range("F14").Value = "Persiapan Juziyah 1 (Juziyah tanggal 13)"

# --- Row 15 (new) ---
$ws.Range("A15").Value = "Idris"
$ws.Range("B15").Value = "Juz 29, 7B"
$ws.Range("C15").Value = "14 halaman"
$ws.Range("D15").Value = "6 kali"
$ws.Range("E15").Value = "2 juz"

# --- Row 16 (new) ---
$ws.Range("A16").Value = "Fahri"
$ws.Range("B16").Value = "Juz 29, 5A"
$ws.Range("C16").Value = "12 halaman"
$ws.Range("D16").Value = "5 kali"
$ws.Range("E16").Value = "2 juz"

# --- Row 17 (new) ---
$ws.Range("A17").Value = "Syahid"
$ws.Range("B17").Value = "Juz 29, 9A"
$ws.Range("C17").Value = "13 halaman"
$ws.Range("D17").Value = "6 kali"
$ws.Range("E17").Value = "2 juz"

# --- Row 18 (new) ---
$ws.Range("A18").Value = "Kamil"
$ws.Range("B18").Value = "Juz 30, 8A"
$ws.Range("C18").Value = "6 halaman"
$ws.Range("D18").Value = "6 kali"
$ws.Range("E18").Value = "1 juz"

# --- Column width for column E widened to fit new (longer) content ---
# Target stored width is 26.42578125 "characters" in the saved file. This
# engine's ColumnWidth setter quantizes to the nearest 1/6 and then stores
# (rounded-value + 5/6), so we back out the input that lands closest to the
# desired stored width.
$ws.Columns.Item(5).ColumnWidth = 25.59244791666667

# --- Selection moves to H15 ---
$ws.Range("H15").Select()
